$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44405
$ws.Range("I2").Value = 'Segunda'
$ws.Range("J2").Value = 140
$ws.Range("D3").Value = 44211
$ws.Range("H3").Value = 'Cultivar XV región'
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 4750
$ws.Range("N3").Value = '$/caja 10 kilos'
$ws.Range("O3").Value = 'Región de Arica y Parinacota'
$ws.Range("P3").Value = 475
$ws.Range("Q3").Value = 10
$ws.Range("D4").Value = 44433
$ws.Range("J4").Value = 100
$ws.Range("D5").Value = 44433
$ws.Range("H5").Value = 'Cultivar IV Región'
$ws.Range("I5").Value = 'Tercera'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = '$/bandeja 18 kilos'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 806
$ws.Range("Q5").Value = 18
$ws.Range("D6").Value = 44398
$ws.Range("M6").Value = 17500
$ws.Range("P6").Value = 972
$ws.Range("D7").Value = 44398
$ws.Range("H7").Value = 'Cultivar IV Región'
$ws.Range("I7").Value = 'Segunda'
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("N7").Value = '$/bandeja 18 kilos'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 861
$ws.Range("Q7").Value = 18
$ws.Range("D8").Value = 44454
$ws.Range("H8").Value = 'Cultivar IV Región'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19500
$ws.Range("N8").Value = '$/bandeja 18 kilos'
$ws.Range("O8").Value = 'Provincia de Limarí'
$ws.Range("P8").Value = 1083
$ws.Range("Q8").Value = 18
$ws.Range("D9").Value = 44221
$ws.Range("J9").Value = 140
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5500
$ws.Range("P9").Value = 550
$ws.Range("D10").Value = 44391
$ws.Range("H10").Value = 'Cultivar IV Región'
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("N10").Value = '$/bandeja 18 kilos'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 861
$ws.Range("Q10").Value = 18
$ws.Range("D11").Value = 44554
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5500
$ws.Range("P11").Value = 550
$ws.Range("D12").Value = 44363
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("P12").Value = 806
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17500
$ws.Range("P13").Value = 972
$ws.Range("D14").Value = 44435
$ws.Range("I14").Value = 'Tercera'
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 806
$ws.Range("D15").Value = 44526
$ws.Range("H15").Value = 'Cultivar XV región'
$ws.Range("I15").Value = 'Primera'
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5500
$ws.Range("M15").Value = 5250
$ws.Range("N15").Value = '$/caja 10 kilos'
$ws.Range("O15").Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 525
$ws.Range("Q15").Value = 10
$ws.Range("D16").Value = 44526
$ws.Range("H16").Value = 'Cultivar XV región'
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = 4250
$ws.Range("N16").Value = '$/caja 10 kilos'
$ws.Range("O16").Value = 'Región de Arica y Parinacota'
$ws.Range("P16").Value = 425
$ws.Range("Q16").Value = 10
$ws.Range("D17").Value = 44526
$ws.Range("H17").Value = 'Cultivar XV región'
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3250
$ws.Range("N17").Value = '$/caja 10 kilos'
$ws.Range("O17").Value = 'Región de Arica y Parinacota'
$ws.Range("P17").Value = 325
$ws.Range("Q17").Value = 10
$ws.Range("D18").Value = 44412
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 972
$ws.Range("D19").Value = 44533
$ws.Range("H19").Value = 'Cultivar XV región'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = 6500
$ws.Range("N19").Value = '$/caja 10 kilos'
$ws.Range("O19").Value = 'Región de Arica y Parinacota'
$ws.Range("P19").Value = 650
$ws.Range("Q19").Value = 10
$ws.Range("D20").Value = 44533
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = 4500
$ws.Range("P20").Value = 450
$ws.Range("D21").Value = 44377
$ws.Range("H21").Value = 'Cultivar IV Región'
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17600
$ws.Range("N21").Value = '$/bandeja 18 kilos'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 978
$ws.Range("Q21").Value = 18
